# REFramework_Porreta/Data/Config.xlsx - "Finish sending BRE notification #15"
#
# Summary of changes applied to the "Settings" sheet:
#  1. Insert a new row 7: ProcessToKill / notepad.exe / description.
#  2. Update the BRE_Notification_BodyHtmlFilePath value to the new path
#     (now under Data\HtmlTemplates\).
#  3. Remove the mailto: hyperlink + its special "Hiperlink" style from the
#     BRE_Notification_Recipients value cell.
#  4. Rename BRE_Notification_Attachment(s)Folders/Files -> singular, and
#     add them as two new rows after BCC (with descriptions).
#  5. Keep new rows styled consistently with their neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- 1. Insert new "ProcessToKill" row at row 7 ------------------------------
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").RowHeight = 14.25
$ws.Range("A7").Value2 = "ProcessToKill"
$ws.Range("B7").Value2 = "notepad.exe"
$ws.Range("C7").Value2 = "Define process to be killed in KillAllProcesses.xaml if in_ProcessToKill is null"

# --- 2. Update the BRE body html path (now row 16 after the insert) --------
$ws.Range("B16").Value2 = "Data\HtmlTemplates\BRE_Body.html"

# --- 3. Remove hyperlink + hyperlink style from the recipients cell --------
# (BRE_Notification_Recipients value is now on row 18 after the insert)
$ws.Hyperlinks.Delete()
$ws.Range("B18").Style = "Normal"
$wb.Styles("Hiperlink").Delete()

# --- 4. Rename Attachment(s)Folders/Files labels & add description rows ----
$ws.Range("A21").Value2 = "BRE_Notification_AttachmentFolders"
$ws.Range("A21").Style = "Normal"
$ws.Range("C21").Value2 = "; separated list of folders whose files will be added as attachments in the business exception email"
$ws.Range("C21").Style = "Normal"

$ws.Range("A22").Value2 = "BRE_Notification_AttachmentFiles"
$ws.Range("A22").Style = "Normal"
$ws.Range("C22").Value2 = "; separated list of files that will be added as attachments in the business exception email"
$ws.Range("C22").Style = "Normal"

# Match the formatting of the other A/C label+description columns in this block
$ws.Range("A21").Font.Size = $ws.Range("A19").Font.Size
$ws.Range("A21").Font.Bold = $ws.Range("A19").Font.Bold
$ws.Range("C21").Font.Size = $ws.Range("C19").Font.Size
$ws.Range("A22").Font.Size = $ws.Range("A19").Font.Size
$ws.Range("C22").Font.Size = $ws.Range("C19").Font.Size

# --- 5. View state: clear the scrolled/top-left cell, move selection -------
[void]$ws.Range("B14").Select()

Write-Output "applied"
